$wb = $excel.ActiveWorkbook

# 1. Rename sheet "Etapa 2" to "Etapa-2"
$wsEtapa2 = $wb.Worksheets.Item("Etapa 2")
$wsEtapa2.Name = "Etapa-2"

# 2. Fix "PUSH A " (trailing space) -> "PUSH A" on that same sheet, cell A90
$wsEtapa2.Range("A90").Value = "PUSH A"

# 3. On "Etapa 1", fill in the missing Operando 2 ("Ins") for the jump
#    instructions in rows 66-71 so the CONCATENATE formula in column A
#    produces "JNE Ins", "JGT Ins", "JLT Ins", "JGE Ins", "JLE Ins", "JCR Ins".
#    Copy the formatting from the neighboring "Instruccion" cell (column C)
#    first so the new cell picks up the same style used by the rest of the
#    "Operando 2" column (e.g. D64/D65), then write the text value.
$wsEtapa1 = $wb.Worksheets.Item("Etapa 1")
foreach ($r in 66..71) {
    $wsEtapa1.Range("C$r").Copy()
    $wsEtapa1.Range("D$r").PasteSpecial(-4122)
    $wsEtapa1.Range("D$r").Value = "Ins"
}
